$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (67) down to
# the two new rows (68-69) so the new cells inherit the same cell styles
# (bold/bordered index column, date-formatted match-time column, etc.)
$ws.Range("A67:V67").Copy()
$ws.Range("A68:V69").PasteSpecial(-4122)

# Row 68: Chabab Mohammedia 0 x 2 Renaissance Zemamra
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "morocco"
$ws.Cells.Item(68, 3).Value = "botola-pro"
$ws.Cells.Item(68, 4).Value = "2023-2024"
$ws.Cells.Item(68, 5).Value = 45242.76041666666
$ws.Cells.Item(68, 6).Value = "Chabab Mohammedia"
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = "Renaissance Zemamra"
$ws.Cells.Item(68, 9).Value = 2
$ws.Cells.Item(68, 10).Value = 2.3
$ws.Cells.Item(68, 11).Value = "11/11/2023 06:42"
$ws.Cells.Item(68, 12).Value = 2.36
$ws.Cells.Item(68, 13).Value = "12/11/2023 18:03"
$ws.Cells.Item(68, 14).Value = 2.83
$ws.Cells.Item(68, 15).Value = "11/11/2023 06:42"
$ws.Cells.Item(68, 16).Value = 2.92
$ws.Cells.Item(68, 17).Value = "12/11/2023 18:03"
$ws.Cells.Item(68, 18).Value = 3.03
$ws.Cells.Item(68, 19).Value = "11/11/2023 06:42"
$ws.Cells.Item(68, 20).Value = 3.28
$ws.Cells.Item(68, 21).Value = "12/11/2023 18:03"
$ws.Cells.Item(68, 22).Value = "https://www.betexplorer.com/football/morocco/botola-pro/chabab-mohammedia-renaissance-zemamra/QFTKNSFk/"

# Row 69: Hassania Agadir 1 x 1 IR Tanger
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "morocco"
$ws.Cells.Item(69, 3).Value = "botola-pro"
$ws.Cells.Item(69, 4).Value = "2023-2024"
$ws.Cells.Item(69, 5).Value = 45242.85416666666
$ws.Cells.Item(69, 6).Value = "Hassania Agadir"
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = "IR Tanger"
$ws.Cells.Item(69, 9).Value = 1
$ws.Cells.Item(69, 10).Value = 2.56
$ws.Cells.Item(69, 11).Value = "11/11/2023 08:43"
$ws.Cells.Item(69, 12).Value = 2.34
$ws.Cells.Item(69, 13).Value = "12/11/2023 20:28"
$ws.Cells.Item(69, 14).Value = 2.8
$ws.Cells.Item(69, 15).Value = "11/11/2023 08:43"
$ws.Cells.Item(69, 16).Value = 2.69
$ws.Cells.Item(69, 17).Value = "12/11/2023 20:28"
$ws.Cells.Item(69, 18).Value = 2.72
$ws.Cells.Item(69, 19).Value = "11/11/2023 08:43"
$ws.Cells.Item(69, 20).Value = 3.69
$ws.Cells.Item(69, 21).Value = "12/11/2023 20:28"
$ws.Cells.Item(69, 22).Value = "https://www.betexplorer.com/football/morocco/botola-pro/hassania-agadir-ir-tanger/8vJPMnVe/"
